$p = $ppt.ActivePresentation

# Remove the trailing slides (4 through the last one), keeping only the
# first three slides (businessoverview, revenue placeholders stay, but
# the "figures" slides and the operatingincome/cashflow slides go away).
for ($i = $p.Slides.Count; $i -ge 4; $i--) {
    $p.Slides.Item($i).Delete()
}

# Slide 2 ("businessoverview") -> update title and body text.
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = '{"body": "Test Item"}'
$s2.Shapes.Item(2).TextFrame.TextRange.Text = 'Please provide the financial documents that contain information about the "Test Item" of the company, so I can analyze how it has changed over time. You can upload the documents here.'

# Slide 3 ("revenue") -> update title and body text.
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = '{"body": "Test Item"}'
$s3.Shapes.Item(2).TextFrame.TextRange.Text = 'It seems that there are no documents currently uploaded. Please upload the financial statements or any relevant documents that contain information about the "Test Item" of the company. Once the documents are uploaded, I can analyze the data and provide a detailed analysis of how the "Test Item" has changed over time.'
